$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Target ordered list of (level, text) pairs for the rebuilt task list.
#   level 0 / 1 -> list indent level (w:ilvl)
#   level -1    -> trailing sentinel paragraph, ListParagraph style but no
#                  numbering at all (mirrors the blank paragraph Word leaves
#                  at the very end of the body)
# ---------------------------------------------------------------------------
$items = @(
    ,@(0, "Research difference between Oracle and MySQL and make a decision between the two.")
    ,@(0, "Think about how the database will be structured")
    ,@(0, "Look into Python modules which may be useful.")
    ,@(1, "Distance calculation")
    ,@(1, "More?")
    ,@(0, "Set up base Python project with GIT version control")
    ,@(0, "Install database software")
    ,@(0, "Investigate similar studies for ideas. Document key findings.")
    ,@(0, "Find the best source for the Premier League statistics required.")
    ,@(0, "Investigate Machine Learning, ML, predictive algorithms and decide on 4 possibilities.")
    ,@(0, "Write TMA01")
    ,@(0, "Write TMA02")
    ,@(0, "Write TMA03")
    ,@(0, "Write EMA")
    ,@(0, "Reflect on progress to date, what went well/bad (Needs to be done multiple times in the project.")
    ,@(0, "Revaluate project after each TMA given feedback from tutor for TMA. Make sure project still makes sense.")
    ,@(0, "CF1 – setting up the database and data")
    ,@(1, "Gather data required.")
    ,@(1, "Cleanse data.")
    ,@(1, "Insert data into database.")
    ,@(0, "CF2 – work on algorithm 1")
    ,@(1, "Plan how algorithm will work.")
    ,@(1, "Do some more research on top of what has already been done.")
    ,@(1, "Code the algorithm")
    ,@(1, "Test the findings")
    ,@(1, "Evaluate")
    ,@(1, "Produce graphs and report to show successfulness of the algorithms ")
    ,@(0, "CF3 – work on algorithm 2")
    ,@(1, "Plan how algorithm will work.")
    ,@(1, "Do some more research on top of what has already been done.")
    ,@(1, "Code the algorithm")
    ,@(1, "Test the findings")
    ,@(1, "Evaluate")
    ,@(1, "Produce graphs and report to show successfulness of the algorithms ")
    ,@(0, "CF4 – work on algorithm 3")
    ,@(1, "Plan how algorithm will work.")
    ,@(1, "Do some more research on top of what has already been done.")
    ,@(1, "Code the algorithm")
    ,@(1, "Test the findings")
    ,@(1, "Evaluate")
    ,@(1, "Produce graphs and report to show successfulness of the algorithms ")
    ,@(0, "CF5 – work on algorithm 4")
    ,@(1, "Plan how algorithm will work.")
    ,@(1, "Do some more research on top of what has already been done.")
    ,@(1, "Code the algorithm")
    ,@(1, "Test the findings")
    ,@(1, "Evaluate")
    ,@(1, "Produce graphs and report to show successfulness of the algorithms ")
    ,@(-1, "")
)

# The first two list paragraphs ("Define the goals..." and "Research SDLC
# choices...") are untouched by the edit, so keep paragraphs 1-2 (the
# heading "Task list" is paragraph 1, "Define the goals..." is paragraph 2,
# "Research SDLC choices..." is paragraph 3) and replace everything from
# paragraph 4 to the end of the document with the rebuilt list above.
$firstToReplace = $d.Paragraphs.Item(4)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$killRange = $d.Range($firstToReplace.Range.Start, $lastPara.Range.End)
$killRange.Delete()

# Grow the list by repeatedly inserting a new paragraph after the current
# last paragraph (which inherits the ListParagraph style + numPr of its
# predecessor), then setting its text/level explicitly.
$cur = $d.Paragraphs.Item($d.Paragraphs.Count)
foreach ($item in $items) {
    $level = $item[0]
    $text = $item[1]

    $cur.Range.InsertParagraphAfter()
    $cur = $d.Paragraphs.Item($d.Paragraphs.Count)
    $cur.Range.Text = $text

    if ($level -eq -1) {
        $cur.Range.ListFormat.RemoveNumbers()
    } else {
        # ListLevelNumber is 1-based (1 -> w:ilvl 0, 2 -> w:ilvl 1, ...),
        # and sets the level absolutely (unlike ListIndent/ListOutdent,
        # which shift relative to the current level).
        $cur.Range.ListFormat.ListLevelNumber = $level + 1
    }
}

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
